# Applies "added new investment validation exceptions":
#  - Removes the "TealBond" investment row (was ID 5)
#  - Removes the "string" investment row (was ID 9) — an invalid placeholder entry
#  - Adds a new "GovBond" investment row (ID 34) right after GlueStock

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "string" row (originally row 10) first, since it is below the
# "TealBond" row and removing it first avoids having to recompute the
# TealBond row index afterwards.
$ws.Rows.Item(10).Delete()

# Remove the "TealBond" row (row 6).
$ws.Rows.Item(6).Delete()

# Insert a new row at position 4 (pushes VegaStock and the rows below it down
# by one) and fill it in with the new GovBond investment data.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = 34
$ws.Cells.Item(4, 2).Value = "GovBond"
$ws.Cells.Item(4, 3).Value = 12
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 12
$ws.Cells.Item(4, 6).Value = 72
